# Change data source to CDC
# - Round several already-present percentage values in rows 664-672 (and one
#   value in row 661) to a small number of decimals (as delivered by the new
#   CDC data source).
# - Append 9 new data rows (673-681) for the next reporting date (3/2/2021).
# - Extend the sheet dimension and the FINAL_AGEGROUPS defined name to the
#   new extent (A1:H681).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FINAL_AGEGROUPS")

# ---------------------------------------------------------------------------
# 1. Update existing cells whose precision was reduced by the new source.
# ---------------------------------------------------------------------------
$updates = @{
    "F661" = 0.3
    "F664" = 0.562
    "H664" = 0.548
    "F665" = 6.826
    "H665" = 8.191
    "F666" = 8.948
    "H666" = 10.81
    "F667" = 9.908
    "H667" = 11.55
    "F668" = 11.05
    "H668" = 12.51
    "F669" = 18.49
    "H669" = 10.93
    "F670" = 29.74
    "H670" = 26.75
    "F671" = 14.44
    "H671" = 18.69
    "F672" = 0.012
    "H672" = 0.002
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---------------------------------------------------------------------------
# 2. Append the new rows (673-681) for DATE 3/2/2021 (serial 44257).
# ---------------------------------------------------------------------------
$newRows = @(
    @("16-20",   44257, 7990,   0.556535742619981, 5261,   0.563728904366461, 2718,   0.545514756796355),
    @("21-30",   44257, 103730, 7.22521308910771,  63217,  6.77385480846504,  40383,  8.10504872101075),
    @("31-40",   44257, 136180, 9.48548653691977,  82958,  8.88915081703724,  52965,  10.6303123965118),
    @("41-50",   44257, 149071, 10.3833967068965,  92093,  9.86798821323332,  56617,  11.3632851308091),
    @("51-60",   44257, 165141, 11.5027370553199,  103092, 11.0465577283686,  61651,  12.3736314463768),
    @("61-70",   44257, 234761, 16.3520509978985,  179048, 19.185427270292,   55177,  11.0742706901223),
    @("71-80",   44257, 412161, 28.7086768728403,  275135, 29.4813822662738,  135587, 27.2129173398629),
    @("81+",     44257, 226525, 15.7783803625771,  132352, 14.1818376640771,  93133,  18.6922096558922),
    @("PENDING", 44257, 108,    0.00752263582014492, 94,  0.0100723278864184, 14,  0.00280986261778844)
)

$startRow = 673
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A: age-group label (plain text, default style).
    $ws.Cells.Item($r, 1).Value = $row[0]

    # Column B: DATE - copy the date number-format from the row above so the
    # new cell uses the same style (rather than Excel creating a new one),
    # then set the actual date value.
    $ws.Cells.Item($r - 1, 2).Copy() | Out-Null
    $ws.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 2).Value = $row[1]

    # Columns C-H: numeric data, default style.
    for ($c = 3; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Update the FINAL_AGEGROUPS defined name to cover the new range.
# ---------------------------------------------------------------------------
$wb.Names.Item("FINAL_AGEGROUPS").RefersTo = "='FINAL_AGEGROUPS'!`$A`$1:`$H`$681"
